$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.991.01'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.88%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.113.14'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.63'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.52'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.44%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.109.08'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.99%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.44'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.60%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.66%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.15'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.22%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.628.50'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.993.63'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.110.14'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.44'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '476.52'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.712'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.90'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +6.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.50'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.85'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.30%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.92'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.43'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.93'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.91%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.67'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.16%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0938'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -7.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.978'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '47.57'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.08'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.92'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.96%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.60'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.796.09'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.75%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '379.61'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.55'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -11.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '136.11'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.02'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.20'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.00%  '
